$wb = $excel.ActiveWorkbook

# Update "展览" sheet (exhibitions) and "全部类型" sheet (all types) which
# both mirror the same underlying rows for this event.
foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 449
    $ws.Range("F3").Value = 15
}
